# Atualizando o arquivo XLSX
# Update several betting-odds values in row 2 of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.67
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.38
$ws.Range("K2").Value = 2.05
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("AD2").Value = 6.5
$ws.Range("AG2").Value = 12
$ws.Range("AH2").Value = 26
$ws.Range("AT2").Value = 2.5
$ws.Range("AX2").Value = 29
